# The workbook is already open; "Repayment schedule" is the active sheet
# (activeTab points at it already), but look it up by name to be safe.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Remember the width of column M (the column that will end up just to the
# left of the newly inserted column) so it can be copied onto the new
# column, mimicking Excel's native "insert column" behaviour.
$leftColWidth = $ws.Columns.Item(13).ColumnWidth

# Insert a new blank column before column N. This shifts the existing
# N/O/P columns (Late / heading / Outstanding) one position to the right,
# turning them into O/P/Q, and leaves the newly created N column blank.
$ws.Range("N1").EntireColumn.Insert() | Out-Null

# Excel carries the column width of the column to the left into the
# newly inserted column - replicate that here.
$ws.Columns.Item(14).ColumnWidth = $leftColWidth

# Update the current selection on the sheet (reflecting where the user
# left the cursor after performing the edit).
$ws.Range("R6").Select() | Out-Null
